$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.609.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.084.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +10.24%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.364"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0723"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0987"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.380.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.829"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.068.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.33%  "
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.675.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "238.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.56%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.50%  "
$ws.Range("E30").Value = "  -5.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +55.17%  "
$ws.Range("E32").Value = "  -1.77%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0584"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0906"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +18.12%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +18.50%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.903"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.35%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.37%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.97%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.328.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0822"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.277.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.09%  "
$ws.Range("E50").Value = "  -5.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +15.22%  "
